# UART Protocol.xlsx - add LCD scroll-buffer commands, rename LCD group to LCDs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UART Commands")
$ws.Activate() | Out-Null

# Rename the "LCD" group label (B24, merged B24:B25) to "LCDs"
$ws.Cells.Item(24, 2).Value = "LCDs"

# Row 26: CREATE_SCROLL_BUFFER (0x60)
$ws.Cells.Item(26, 3).Value = "0x60"
$ws.Cells.Item(26, 4).Value = "CREATE_SCROLL_BUFFER"
$ws.Cells.Item(26, 5).Value = 8
$ws.Cells.Item(26, 6).Value = 1
$ws.Cells.Item(26, 7).Value = "Create a scroll buffer.`nParameters:`nByte 0 - X position high byte`nByte 1 - X position low byte`nByte 2 - Y position high byte`nByte 3 - Y position low byte`nByte 4 - Width high byte`nByte 5 - Width low byte`nByte 6 - Height high byte`nByte 7 - Height low byte`nResponse:`nByte 0 - The ID of the created buffer or 0 on failure"
$ws.Rows.Item(26).RowHeight = 180

# Row 27: DELETE_SCROLL_BUFFER (0x61)
$ws.Cells.Item(27, 3).Value = "0x61"
$ws.Cells.Item(27, 4).Value = "DELETE_SCROLL_BUFFER"
$ws.Cells.Item(27, 5).Value = 1
$ws.Cells.Item(27, 6).Value = 1
$ws.Cells.Item(27, 7).Value = "Delete a scroll buffer.`nParameters:`nByte 0 - The ID of the buffer to be deleted`nResponse:`nByte 0 - 1 on success, 0 on failure"
$ws.Rows.Item(27).RowHeight = 75

# Row 28: SET_DESTINATION_BUFFER (0x62)
$ws.Cells.Item(28, 3).Value = "0x62"
$ws.Cells.Item(28, 4).Value = "SET_DESTINATION_BUFFER"
$ws.Cells.Item(28, 5).Value = 1
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(28, 7).Value = "Set the SPI destination buffer.`nParameters:`nByte 0 - The ID of the destination buffer`nResponse:`nByte 0 - 1 on success, 0 on failure"
$ws.Rows.Item(28).RowHeight = 75

# Update the view: scroll so row 25 is at the top, and select C29 (the next empty row)
$ws.Range("C29").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
